$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'29.959.71"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.22%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'1.877.52"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -0.91%  "
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.7408"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -4.22%  "
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'242.27"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  -0.86%  "
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'1.000"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  -0.05%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.3155"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +0.86%  "
$c.Style = "Normal"

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c = $ws.Range("D9")
$c.Value = "'0.07176"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -0.47%  "
$c.Style = "Normal"

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c = $ws.Range("D10")
$c.Value = "'24.81"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  -3.46%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.08433"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -4.93%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.7558"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -2.11%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'5.427"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -0.01%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'1.867.81"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -5.41%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'92.92"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -1.61%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'29.957.47"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +0.00%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'6.106"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  -1.35%  "
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'  -2.15%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'244.07"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -0.62%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.000007845"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -0.26%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'0.9994"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -0.13%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'2.114.80"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -3.88%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'7.998"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -2.01%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'0.9973"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -0.38%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'0.1566"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  -2.23%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'9.332"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -2.00%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'164.78"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  +1.34%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'18.70"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -0.63%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'2.043"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -0.01%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'1.478"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +3.52%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'4.626"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +1.34%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'1.532"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -0.74%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'4.288"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +4.26%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'0.05334"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  -2.90%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'1.241"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -0.66%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'0.7599"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  +1.17%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'0.9985"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'2.696"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  -0.74%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.01956"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'2.751"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -1.21%  "
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'  -0.24%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'1.110.10"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  +1.47%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'6.144"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +1.90%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'72.65"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  -1.82%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.8637"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  +1.05%  "
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'  +0.08%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'103.33"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  +0.68%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'7.713"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  +1.20%  "
$c.Style = "Normal"

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D49")
$c.Value = "'1.851"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -1.92%  "
$c.Style = "Normal"

$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$c = $ws.Range("D50")
$c.Value = "'3.077"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +4.01%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'2.014.32"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -2.09%  "
$c.Style = "Normal"
